$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# --- 1. Row 58: mark the map-position task as Finished instead of "In the work" ---
$ws.Range("E58").Value2 = "Finished"

# --- 2. Insert two new rows before row 60 (mirrors a user selecting rows 60:61 and
#        choosing Insert; Excel shifts everything below down by two and auto-adjusts
#        merged cells / SUM ranges / defined names that touch the insertion point). ---
$ws.Rows("60:61").Insert()

# --- 3. Row 59 already existed (blank) and keeps almost all of its formatting; only
#        its time cell needs the "time" number format used by the rest of the block. ---
$ws.Range("F58").Copy()
$ws.Range("F59").PasteSpecial(-4122)

# --- 4. The two freshly inserted rows (60 and 61) need the same direct-format styles
#        used by the rest of this data block. Build them column by column from the
#        matching cells still present on the sheet. ---
$ws.Range("A59").Copy()
$ws.Range("A60:A61").PasteSpecial(-4122)

$ws.Range("B59").Copy()
$ws.Range("B60:B61").PasteSpecial(-4122)

$ws.Range("C62").Copy()
$ws.Range("C60:C61").PasteSpecial(-4122)

$ws.Range("D59").Copy()
$ws.Range("D60:D61").PasteSpecial(-4122)

$ws.Range("E62").Copy()
$ws.Range("E60:E61").PasteSpecial(-4122)

$ws.Range("F58").Copy()
$ws.Range("F60:F61").PasteSpecial(-4122)

$ws.Range("G59").Copy()
$ws.Range("G60:G61").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- 5. Row 62 (the old end-of-block row, now shifted down) keeps its thick-bottom
#        border but its first cell should use the plain interior style like the rest
#        of column A in this block. ---
$ws.Range("A59").Copy()
$ws.Range("A62").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 6. Fill in the three entries describing the Issue #7 work (text entered in the
#        same order as the author - row 61 first, then row 59, then row 60). ---
$ws.Range("A61").Value2 = "Coding"
$ws.Range("C61").Value2 = 35
$ws.Range("D61").Value2 = "Implementation of Issue #7 to the program"
$ws.Range("E61").Value2 = "Finished"
$ws.Range("F61").Value2 = 0.69097222222222221

$ws.Range("A59").Value2 = "Coding"
$ws.Range("C59").Value2 = 30
$ws.Range("D59").Value2 = "Creation of new enemy types for issue #7"
$ws.Range("E59").Value2 = "Finished"
$ws.Range("F59").Value2 = 0.64583333333333337

$ws.Range("A60").Value2 = "Coding"
$ws.Range("C60").Value2 = 30
$ws.Range("D60").Value2 = "Fixed small issues with map position"
$ws.Range("E60").Value2 = "Finished"
$ws.Range("F60").Value2 = 0.66666666666666663

# --- 7. Update the print area to match the new sheet extent (last row moved from 68
#        to 70 because of the two inserted rows). ---
$ws.PageSetup.PrintArea = "`$A`$1:`$G`$70"

# --- 8. Restore the active selection recorded for this sheet after the edit. ---
$ws.Range("E62").Select()
